$wb = $excel.ActiveWorkbook

# Sheet "展览" - update column F (想去人数) values
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 70
$ws1.Range("F4").Value = 232
$ws1.Range("F6").Value = 9934
$ws1.Range("F7").Value = 899
$ws1.Range("F9").Value = 1240
$ws1.Range("F10").Value = 4474
$ws1.Range("F15").Value = 59
$ws1.Range("F18").Value = 570
$ws1.Range("F19").Value = 111
$ws1.Range("F21").Value = 1489

# Sheet "全部类型" - update column F (想去人数) values
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 70
$ws4.Range("F5").Value = 232
$ws4.Range("F7").Value = 9934
$ws4.Range("F8").Value = 899
$ws4.Range("F10").Value = 1240
$ws4.Range("F11").Value = 4474
$ws4.Range("F16").Value = 59
$ws4.Range("F19").Value = 570
$ws4.Range("F20").Value = 111
$ws4.Range("F22").Value = 1489
